$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: cohort_year=2023, period_index=3 -> num_customers 60 -> 61
# retention_rate recalculated as num_customers / cohort_size (D16 = 193)
$ws.Range("C16").Value = 61
$ws.Range("E16").Value = 61/193

# Row 22: cohort_year=2025, period_index=0 -> num_customers 96 -> 97
# cohort_size (D22) also updated to match 96 -> 97
$ws.Range("C22").Value = 97
$ws.Range("D22").Value = 97
